$d = $word.ActiveDocument

# Word's Font.Color is a BGR-packed long (wdColor), so convert the
# "345A8A" RGB hex used by the style from the source RGB triplet.
$rgbHex = "345A8A"
$r = [Convert]::ToInt32($rgbHex.Substring(0,2), 16)
$g = [Convert]::ToInt32($rgbHex.Substring(2,2), 16)
$b = [Convert]::ToInt32($rgbHex.Substring(4,2), 16)
$abstractTitleColor = $b * 65536 + $g * 256 + $r

# ------------------------------------------------------------------
# 1. Add the new "Abstract Title" paragraph style, based on Normal,
#    followed by the Abstract style.
# ------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = $abstractTitleColor

# ------------------------------------------------------------------
# 2. Abstract style: reduce space-before from 300 (15pt) to 100 (5pt).
# ------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ------------------------------------------------------------------
# 3. Add the new "Footnote Block Text" paragraph style, based on
#    Footnote Text, followed by Footnote Text again.
# ------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "Styles updated."
